$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "33.10", "1.00").
# Force text format before writing so Excel does not silently convert them
# to Number cells / strip trailing zeros, then clear the format override so
# the cell style index matches the original (unstyled) cells.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "42.203.07"
$ws.Range("E2").Value = "  -1.39%  "
Set-TextValue $ws.Range("D3") "2.272.17"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "299.58"
$ws.Range("E5").Value = "  -1.26%  "
Set-TextValue $ws.Range("D6") "95.69"
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("E7").Value = "  -2.45%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.89%  "
Set-TextValue $ws.Range("D10") "33.10"
$ws.Range("E10").Value = "  -5.12%  "
$ws.Range("E11").Value = "  -0.72%  "
Set-TextValue $ws.Range("D12") "48.22"
$ws.Range("E12").Value = "  -6.93%  "
Set-TextValue $ws.Range("D13") "0.115"
$ws.Range("E13").Value = "  +1.61%  "
Set-TextValue $ws.Range("D14") "16.02"
$ws.Range("E14").Value = "  +1.57%  "
Set-TextValue $ws.Range("D15") "6.67"
$ws.Range("E15").Value = "  -1.06%  "
Set-TextValue $ws.Range("D16") "2.625.54"
$ws.Range("E16").Value = "  -1.62%  "
Set-TextValue $ws.Range("D17") "2.278.33"
$ws.Range("E17").Value = "  -0.98%  "
Set-TextValue $ws.Range("D18") "0.786"
$ws.Range("E18").Value = "  -2.85%  "
Set-TextValue $ws.Range("D19") "42.157.20"
$ws.Range("E19").Value = "  -1.31%  "
Set-TextValue $ws.Range("D20") "11.72"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("E21").Value = "  -1.83%  "
Set-TextValue $ws.Range("D22") "5.98"
$ws.Range("E22").Value = "  -1.70%  "
Set-TextValue $ws.Range("D23") "66.31"
$ws.Range("E23").Value = "  -2.53%  "
Set-TextValue $ws.Range("D24") "235.32"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("E26").Value = "  +0.00%  "
Set-TextValue $ws.Range("D27") "2.45"
$ws.Range("E27").Value = "  -2.82%  "
Set-TextValue $ws.Range("D28") "23.73"
$ws.Range("E28").Value = "  -5.12%  "
Set-TextValue $ws.Range("D29") "167.85"
$ws.Range("E30").Value = "  -4.72%  "
Set-TextValue $ws.Range("D31") "33.59"
$ws.Range("E31").Value = "  -3.54%  "
Set-TextValue $ws.Range("D32") "9.12"
$ws.Range("E32").Value = "  -0.52%  "
Set-TextValue $ws.Range("D33") "1.00"
$ws.Range("E33").Value = "  +0.01%  "
Set-TextValue $ws.Range("D34") "4.69"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("E35").Value = "  -2.50%  "
Set-TextValue $ws.Range("D36") "16.70"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("E39").Value = "  -3.31%  "
Set-TextValue $ws.Range("D40") "0.0987"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("E42").Value = "  -4.55%  "
Set-TextValue $ws.Range("D43") "2.32"
$ws.Range("E43").Value = "  -7.16%  "
Set-TextValue $ws.Range("D44") "1.959.63"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("E45").Value = "  -1.33%  "
Set-TextValue $ws.Range("D46") "17.66"
$ws.Range("E46").Value = "  -5.09%  "
Set-TextValue $ws.Range("D47") "9.59"
$ws.Range("E47").Value = "  -6.31%  "
Set-TextValue $ws.Range("D48") "2.77"
$ws.Range("E48").Value = "  -4.55%  "
$ws.Range("E49").Value = "  -1.49%  "
Set-TextValue $ws.Range("D50") "52.14"
$ws.Range("E50").Value = "  -6.99%  "
$ws.Range("E51").Value = "  -3.61%  "
